# Generate Report for Handoff
# Adds a new tracked file (efd47859-7f9b-474f-b234-27474b71f346.md) to the
# localization-status workbook: one new row on each of the three sheets
# (Overview, zh-cn, de-de), new hyperlinks pointing at the new markdown
# file on GitHub, and the corresponding table/autofilter range + sheet
# dimension growth that comes from extending each table by one row.

$wb = $excel.ActiveWorkbook

$commit = "75417ec7c3cf1676dff5f3beb730110c08a1f40e"
$fileId = "efd47859-7f9b-474f-b234-27474b71f346"
$mdName = "$fileId.md"
$mdDisplay = "e2e\$fileId.md"
$mdUrl = "https://github.com/OpenLocalizationTestOrg/oltest/blob/$commit/e2e/$mdName"

$dateFmt = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------
# Sheet "Overview" (table3 / A1:G2 -> A1:G3)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A3").Value = $mdName
$wsOverview.Range("B3").Value = $mdDisplay
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("D3").Value = ""
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-13 06:49:22"
$wsOverview.Range("G3").NumberFormat = $dateFmt

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdDisplay) | Out-Null
$wsOverview.Range("B3").Style = "HyperLink"

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G3"))

# ---------------------------------------------------------------------
# Sheet "zh-cn" (table1 / A1:P2 -> A1:P3)
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A3").Value = $mdName
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "'False"
$wsZh.Range("G3").Value = "$fileId.6194d4792ff707b256e9b9dfaa8ad62699adf3d4.zh-cn.xlf"
$wsZh.Range("H3").Value = "2016-08-13 06:49:14"
$wsZh.Range("H3").NumberFormat = $dateFmt
$wsZh.Range("I3").Value = ""
$wsZh.Range("J3").Value = ""
$wsZh.Range("K3").Value = "0001-01-01 00:00:00"
$wsZh.Range("K3").NumberFormat = $dateFmt
$wsZh.Range("L3").Value = ""
$wsZh.Range("M3").Value = "'True"
$wsZh.Range("N3").Value = ""
$wsZh.Range("O3").Value = "'False"
$wsZh.Range("P3").Value = ""

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdName) | Out-Null
$wsZh.Range("A3").Style = "HyperLink"

$loZh = $wsZh.ListObjects.Item(1)
$loZh.Resize($wsZh.Range("A1:P3"))

# ---------------------------------------------------------------------
# Sheet "de-de" (table2 / A1:P2 -> A1:P3)
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A3").Value = $mdName
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "'False"
$wsDe.Range("G3").Value = "$fileId.6194d4792ff707b256e9b9dfaa8ad62699adf3d4.de-de.xlf"
$wsDe.Range("H3").Value = "2016-08-13 06:49:22"
$wsDe.Range("H3").NumberFormat = $dateFmt
$wsDe.Range("I3").Value = ""
$wsDe.Range("J3").Value = ""
$wsDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDe.Range("K3").NumberFormat = $dateFmt
$wsDe.Range("L3").Value = ""
$wsDe.Range("M3").Value = "'True"
$wsDe.Range("N3").Value = ""
$wsDe.Range("O3").Value = "'False"
$wsDe.Range("P3").Value = ""

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdName) | Out-Null
$wsDe.Range("A3").Style = "HyperLink"

$loDe = $wsDe.ListObjects.Item(1)
$loDe.Resize($wsDe.Range("A1:P3"))
